$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "RE82"
$ws.Range("B6").Value = "Diana Razo"
$ws.Range("C6").Value = "High ovality on this setup, decreased the water line from 100 to 85. Shrink is still passing with the change."
$ws.Range("D6").Value = "2025-10-20 09:39:35"
